$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7663859724998474
$ws.Range("B1").Value = 2.652182102203369
$ws.Range("C1").Value = 4.884721279144287
$ws.Range("D1").Value = 2.829785823822021
$ws.Range("E1").Value = 1.005312919616699
